# A new "September" entry ("axis") came in on the "2024" sheet.
# This pushes every existing record in row 30 downward (rows 30-65 -> 31-66)
# and the brand new entry is written into the freed-up row 30.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new blank row at row 30; everything below (through the old last
# row 65) shifts down by one row, turning the old dimension A1:Y65 into
# A1:Y66 automatically.
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new September record.
$ws.Cells.Item(30, 18).Value = "axis"
$ws.Cells.Item(30, 19).Value = "2024-09-05 15:57:15"
